# "limit the result of addtion to 10"
#
# 1) Add a new "The test tablet" Heading1 section at the very top of the
#    document describing the physical test device.
# 2) Drop the stray _GoBack bookmark that previously sat at the end of the
#    "Drag and drop" paragraph.
# 3) Append a new "Add learning" Heading1 section (with a "Heading2" sub-
#    section + a table of number-of-objects -> number-of-columns) at the end
#    of the document; the _GoBack bookmark is recreated around the final
#    table cell, which is now the true end-of-document insertion point.

$d = $word.ActiveDocument

# --- 1) Remove the now-obsolete _GoBack bookmark ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2) Insert "The test tablet" section at the very start of the document ---
$section1Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>The test tablet</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr><w:r><w:t>Samsung</w:t></w:r><w:r><w:t xml:space="preserve"> tablet</w:t></w:r><w:r><w:t xml:space="preserve"> SM-T520</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr><w:r><w:t>2560 x 1600 Pixel</w:t></w:r><w:r><w:t xml:space="preserve"> = (1369 * 856dp)</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr><w:r><w:t xml:space="preserve">298.9 </w:t></w:r><w:r><w:t>ppi</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr><w:r><w:t>10.1 inch</w:t></w:r><w:r><w:t>es</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr><w:r><w:t>1px = 0.535dp</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$startRange = $d.Range(0, 0)
[void]$startRange.InsertXML($section1Xml)

# --- 3) Append "Add learning" section (incl. table + relocated bookmark) ---
$section2Xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Add learning</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr><w:r><w:t>Width for operators and tool bar:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr><w:r><w:t>180px(96.3dp)</w:t></w:r><w:r><w:t>. thus width for operands: (1369-96.3*3)/3=360.03</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>The number of columns of the gridview containing objects</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="MediumShading1-Accent6"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4261"/><w:gridCol w:w="4261"/></w:tblGrid><w:tr><w:trPr><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>The number value</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="100000000000" w:firstRow="1" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>Columns</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>1</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>2</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>2</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>4</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>2</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>5</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>6</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>7</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>8</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:t>9</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000100000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="1" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:trPr><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:trPr><w:tc><w:tcPr><w:cnfStyle w:val="001000000000" w:firstRow="0" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="0" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>10</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="4261" w:type="dxa"/></w:tcPr><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLineChars="0" w:firstLine="0"/><w:cnfStyle w:val="000000010000" w:firstRow="0" w:lastRow="0" w:firstColumn="0" w:lastColumn="0" w:oddVBand="0" w:evenVBand="0" w:oddHBand="0" w:evenHBand="1" w:firstRowFirstColumn="0" w:firstRowLastColumn="0" w:lastRowFirstColumn="0" w:lastRowLastColumn="0"/></w:pPr><w:r><w:t>4</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:pStyle w:val="NormalIndent"/><w:ind w:firstLine="600"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$endPos = $d.Content.End
$endRange = $d.Range($endPos, $endPos)
[void]$endRange.InsertXML($section2Xml)

